# Add an "Equipment Type" column to the report, inserted between the
# existing "System #" (D) and "{year}" (E) columns, matching the look
# of the neighbouring header cell (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember column D's width so the newly inserted column E can match it
# (Excel normally carries the left-neighbour's formatting/width into a
# freshly inserted column).
$colDWidth = $ws.Columns("D:D").ColumnWidth

# Insert a new column at E - this shifts the former E:K columns to F:L
# and copies formatting (styles) from column D into the new column E.
$ws.Columns("E:E").Insert()

# Match the new column's width to column D's.
$ws.Columns("E:E").ColumnWidth = $colDWidth

# Header text for the new column (row 2 holds the column headers).
$ws.Range("E2").Value = "Equipment Type"

# Leave the new header cell selected, like the original edit did.
$null = $ws.Range("E2").Select()
